$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text format first, otherwise Excel auto-converts them (losing trailing
# zeros / introducing floating-point rounding) instead of keeping the
# original text-like representation.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated price / volume values
$ws.Range("D2").Value = "52.020.09"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.811.85"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "361.09"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "110.83"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  +2.59%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "40.41"
$ws.Range("E10").Value = "  -5.10%  "
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "19.76"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "3.255.77"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.828.46"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "0.921"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("D18").Value = "51.950.03"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").Value = "3.11"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").Value = "13.29"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").Value = "0.0₃0991"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "271.14"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "69.98"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "26.65"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "0.142"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").Value = "52.11"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").Value = "34.47"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  +10.31%  "
$ws.Range("D36").Value = "0.0845"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").Value = "2.02"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("D40").Value = "18.11"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "0.116"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").Value = "124.80"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "2.26"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").Value = "22.42"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").Value = "2.072.01"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "3.28"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "5.83"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").Value = "0.950"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "9.15"
$ws.Range("E51").Value = "  +0.09%  "
